$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $cellRef, $value)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

    Set-TextValue $ws "D2" "24.744.84"
    Set-TextValue $ws "E2" "  -0.18%  "
    Set-TextValue $ws "D3" "1.675.87"
    Set-TextValue $ws "E3" "  -0.57%  "
    Set-TextValue $ws "D4" "1.003"
    Set-TextValue $ws "E4" "  +0.25%  "
    Set-TextValue $ws "D5" "308.60"
    Set-TextValue $ws "E5" "  +0.74%  "
    Set-TextValue $ws "D6" "0.9968"
    Set-TextValue $ws "E6" "  +0.06%  "
    Set-TextValue $ws "D7" "0.3678"
    Set-TextValue $ws "E7" "  +0.00%  "
    Set-TextValue $ws "D8" "48.02"
    Set-TextValue $ws "E8" "  -4.42%  "
    Set-TextValue $ws "D9" "0.3364"
    Set-TextValue $ws "E9" "  -1.61%  "
    Set-TextValue $ws "D10" "1.174"
    Set-TextValue $ws "E10" "  +1.14%  "
    Set-TextValue $ws "D11" "0.07328"
    Set-TextValue $ws "E11" "  +1.58%  "
    Set-TextValue $ws "D12" "0.9997"
    Set-TextValue $ws "E12" "  +0.24%  "
    Set-TextValue $ws "B13" "Solana"
    Set-TextValue $ws "C13" "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
    Set-TextValue $ws "D13" "20.51"
    Set-TextValue $ws "E13" "  +1.85%  "
    Set-TextValue $ws "B14" "Polkadot"
    Set-TextValue $ws "C14" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
    Set-TextValue $ws "D14" "6.164"
    Set-TextValue $ws "E14" "  +1.05%  "
    Set-TextValue $ws "D15" "6.796"
    Set-TextValue $ws "E15" "  +1.85%  "
    Set-TextValue $ws "D16" "1.676.54"
    Set-TextValue $ws "E16" "  -0.36%  "
    Set-TextValue $ws "D17" "0.00001100"
    Set-TextValue $ws "E17" "  -0.19%  "
    Set-TextValue $ws "D18" "0.06611"
    Set-TextValue $ws "E18" "  -0.70%  "
    Set-TextValue $ws "D19" "0.9967"
    Set-TextValue $ws "E19" "  +0.10%  "
    Set-TextValue $ws "D20" "81.75"
    Set-TextValue $ws "E20" "  +1.01%  "
    Set-TextValue $ws "D21" "16.81"
    Set-TextValue $ws "E21" "  +3.07%  "
    Set-TextValue $ws "D22" "6.182"
    Set-TextValue $ws "E22" "  +2.02%  "
    Set-TextValue $ws "D23" "12.62"
    Set-TextValue $ws "E23" "  +4.58%  "
    Set-TextValue $ws "D24" "24.732.04"
    Set-TextValue $ws "E24" "  +0.16%  "
    Set-TextValue $ws "D25" "2.425"
    Set-TextValue $ws "E25" "  +0.67%  "
    Set-TextValue $ws "D26" "2.704"
    Set-TextValue $ws "E26" "  +2.01%  "
    Set-TextValue $ws "D27" "19.80"
    Set-TextValue $ws "D28" "148.41"
    Set-TextValue $ws "E28" "  -2.70%  "
    Set-TextValue $ws "D29" "129.90"
    Set-TextValue $ws "E29" "  +2.20%  "
    Set-TextValue $ws "D30" "1.864.36"
    Set-TextValue $ws "E30" "  -0.24%  "
    Set-TextValue $ws "D31" "1.222"
    Set-TextValue $ws "E31" "  +24.61%  "
    Set-TextValue $ws "D32" "6.518"
    Set-TextValue $ws "E32" "  +4.78%  "
    Set-TextValue $ws "D33" "4.152"
    Set-TextValue $ws "E33" "  +3.18%  "
    Set-TextValue $ws "B34" "Aptos"
    Set-TextValue $ws "C34" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
    Set-TextValue $ws "D34" "13.40"
    Set-TextValue $ws "E34" "  +8.88%  "
    Set-TextValue $ws "B35" "Stellar"
    Set-TextValue $ws "C35" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
    Set-TextValue $ws "D35" "0.08591"
    Set-TextValue $ws "E35" "  +2.50%  "
    Set-TextValue $ws "D36" "1.728"
    Set-TextValue $ws "E36" "  +2.27%  "
    Set-TextValue $ws "D37" "5.419"
    Set-TextValue $ws "E37" "  +2.54%  "
    Set-TextValue $ws "D38" "0.06440"
    Set-TextValue $ws "E38" "  +1.56%  "
    Set-TextValue $ws "D39" "8.799"
    Set-TextValue $ws "E39" "  +2.25%  "
    Set-TextValue $ws "D40" "0.02340"
    Set-TextValue $ws "E40" "  +1.76%  "
    Set-TextValue $ws "D41" "0.2152"
    Set-TextValue $ws "E41" "  +3.36%  "
    Set-TextValue $ws "D42" "1.233"
    Set-TextValue $ws "E42" "  -0.19%  "
    Set-TextValue $ws "D43" "0.6253"
    Set-TextValue $ws "E43" "  +2.84%  "
    Set-TextValue $ws "D44" "0.9972"
    Set-TextValue $ws "E44" "  +0.17%  "
    Set-TextValue $ws "D45" "13.45"
    Set-TextValue $ws "E45" "  +3.16%  "
    Set-TextValue $ws "D46" "3.784"
    Set-TextValue $ws "E46" "  +0.66%  "
    Set-TextValue $ws "D47" "0.5941"
    Set-TextValue $ws "E47" "  +1.28%  "
    Set-TextValue $ws "D48" "2.041"
    Set-TextValue $ws "E48" "  +2.52%  "
    Set-TextValue $ws "D49" "125.56"
    Set-TextValue $ws "E49" "  +0.34%  "
    Set-TextValue $ws "D50" "0.07163"
    Set-TextValue $ws "E50" "  -0.95%  "
    Set-TextValue $ws "D51" "76.97"
    Set-TextValue $ws "E51" "  +1.97%  "

Write-Host "Applied cryptos price/volume update."
